$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 354 ("spatial scale" cl:10348). This shifts all following rows
# up by one, matching the new dimension A1:CZ396.
$ws.Rows.Item(354).Delete()
